$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: move the "2h21m" entry from column B to column C (retyped as "2h 27m"),
#     keep D13 as-is, and add a new H13 value ---
$ws.Range("B13").ClearContents()
$ws.Range("C13").Value = "2h 27m"
$ws.Range("H13").Value = 7.8

# --- Row 14: add H14 ---
$ws.Range("H14").Value = 7.5

# --- Row 15: add D15 and H15 ---
$ws.Range("D15").Value = "1h 3m"
$ws.Range("H15").Value = 0

# --- Row 16: add C16 and D16 ---
$ws.Range("C16").Value = "1h 56m"
$ws.Range("D16").Value = "1h 19m"

# --- View state: scroll the frozen sheet down and move the active selection ---
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("N14").Select()
